# Auto-generated edit script: refreshes market-price derived columns (H-N)
# on the Leve Profit tables across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values come from a scheduled market-data refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2108.1428
$ws.Range("I40").Value = 1984.3572
$ws.Range("J40").Value = 2355.7144
$ws.Range("K40").Value = 1984.3572
$ws.Range("L40").Value = 2355.7144
$ws.Range("M40").Value = -1809.3572
$ws.Range("N40").Value = -2705.7144
$ws.Range("H49").Value = 200
$ws.Range("I49").Value = 200
$ws.Range("K49").Value = 600
$ws.Range("M49").Value = -464
$ws.Range("H52").Value = 1599.6666
$ws.Range("J52").Value = 3000
$ws.Range("L52").Value = 9000
$ws.Range("N52").Value = -9320
$ws.Range("H64").Value = 3871.4285
$ws.Range("I64").Value = 3040
$ws.Range("J64").Value = 4333.3335
$ws.Range("K64").Value = 3040
$ws.Range("L64").Value = 4333.3335
$ws.Range("M64").Value = -2792
$ws.Range("N64").Value = -4829.3335
$ws.Range("H67").Value = 3871.4285
$ws.Range("I67").Value = 3040
$ws.Range("J67").Value = 4333.3335
$ws.Range("K67").Value = 3040
$ws.Range("L67").Value = 4333.3335
$ws.Range("M67").Value = -2182
$ws.Range("N67").Value = -6049.3335
$ws.Range("H76").Value = 2982.4736
$ws.Range("I76").Value = 2711.8
$ws.Range("K76").Value = 2711.8
$ws.Range("M76").Value = -2396.8
$ws.Range("H79").Value = 2982.4736
$ws.Range("I79").Value = 2711.8
$ws.Range("K79").Value = 2711.8
$ws.Range("M79").Value = -1619.8
$ws.Range("H137").Value = 17243036
$ws.Range("I137").Value = 26316702
$ws.Range("J137").Value = 3070.3
$ws.Range("K137").Value = 78950106
$ws.Range("L137").Value = 9210.900000000001
$ws.Range("M137").Value = -78947556
$ws.Range("N137").Value = -14310.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2150.5833
$ws.Range("I88").Value = 1600
$ws.Range("J88").Value = 2260.7
$ws.Range("K88").Value = 1600
$ws.Range("L88").Value = 2260.7
$ws.Range("M88").Value = -1194
$ws.Range("N88").Value = -3072.7
$ws.Range("H91").Value = 2150.5833
$ws.Range("I91").Value = 1600
$ws.Range("J91").Value = 2260.7
$ws.Range("K91").Value = 1600
$ws.Range("L91").Value = 2260.7
$ws.Range("M91").Value = -196
$ws.Range("N91").Value = -5068.7
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 30000
$ws.Range("J60").Value = 30000
$ws.Range("L60").Value = 30000
$ws.Range("N60").Value = -31198
$ws.Range("H105").Value = 2962.3333
$ws.Range("I105").Value = 2450
$ws.Range("J105").Value = 3218.5
$ws.Range("K105").Value = 2450
$ws.Range("L105").Value = 3218.5
$ws.Range("M105").Value = -703
$ws.Range("N105").Value = -6712.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 8000
$ws.Range("J23").Value = 8000
$ws.Range("L23").Value = 8000
$ws.Range("N23").Value = -8480
$ws.Range("H27").Value = 8000
$ws.Range("J27").Value = 8000
$ws.Range("L27").Value = 8000
$ws.Range("N27").Value = -8384
$ws.Range("H56").Value = 12551.5
$ws.Range("I56").Value = 5000
$ws.Range("J56").Value = 20103
$ws.Range("K56").Value = 5000
$ws.Range("L56").Value = 20103
$ws.Range("M56").Value = -4155
$ws.Range("N56").Value = -21793
$ws.Range("H62").Value = 500000
$ws.Range("I62").Value = 500000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 500000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -499376
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H65").Value = 500000
$ws.Range("I65").Value = 500000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 2500000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -2496880
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H122").Value = 5564.654
$ws.Range("I122").Value = 6741.684
$ws.Range("J122").Value = 2369.8572
$ws.Range("K122").Value = 20225.052
$ws.Range("L122").Value = 7109.571599999999
$ws.Range("M122").Value = -17775.052
$ws.Range("N122").Value = -12009.5716
$ws.Range("N62").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2793.647
$ws.Range("I70").Value = 1249.2
$ws.Range("K70").Value = 3747.6
$ws.Range("M70").Value = -3432.6
$ws.Range("H73").Value = 2793.647
$ws.Range("I73").Value = 1249.2
$ws.Range("K73").Value = 3747.6
$ws.Range("M73").Value = -2655.6
$ws.Range("H93").Value = 7520
$ws.Range("I93").Value = 1500
$ws.Range("J93").Value = 9025
$ws.Range("K93").Value = 4500
$ws.Range("L93").Value = 27075
$ws.Range("M93").Value = -2628
$ws.Range("N93").Value = -30819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 19987.5
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 19987.5
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 19987.5
$ws.Range("N58").Value = -20541.5
$ws.Range("H70").Value = 4393
$ws.Range("I70").Value = 4481.1665
$ws.Range("J70").Value = 4304.8335
$ws.Range("K70").Value = 4481.1665
$ws.Range("L70").Value = 4304.8335
$ws.Range("M70").Value = -4211.1665
$ws.Range("N70").Value = -4844.8335
$ws.Range("H73").Value = 4393
$ws.Range("I73").Value = 4481.1665
$ws.Range("J73").Value = 4304.8335
$ws.Range("K73").Value = 4481.1665
$ws.Range("L73").Value = 4304.8335
$ws.Range("M73").Value = -3545.1665
$ws.Range("N73").Value = -6176.8335
$ws.Range("H80").Value = 202121.2
$ws.Range("J80").Value = 252076.5
$ws.Range("L80").Value = 252076.5
$ws.Range("N80").Value = -254072.5
$ws.Range("H82").Value = 34750
$ws.Range("J82").Value = 40000
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40766
$ws.Range("H83").Value = 202121.2
$ws.Range("J83").Value = 252076.5
$ws.Range("L83").Value = 1260382.5
$ws.Range("N83").Value = -1270366.5
$ws.Range("H85").Value = 34750
$ws.Range("J85").Value = 40000
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42652
$ws.Range("M58").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 35000
$ws.Range("J82").Value = 35000
$ws.Range("L82").Value = 35000
$ws.Range("N82").Value = -35766
$ws.Range("H85").Value = 35000
$ws.Range("J85").Value = 35000
$ws.Range("L85").Value = 35000
$ws.Range("N85").Value = -37652
$ws.Range("H95").Value = 41211
$ws.Range("J95").Value = 41211
$ws.Range("L95").Value = 41211
$ws.Range("N95").Value = -46703
